# Update curva_forward worksheet: shift batch dates/ids forward one day and refresh vertexValue (F) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 69
$ws.Cells.Item(2, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(2, 8).Value = "2023-01-11T05:32:11.915Z"
$ws.Cells.Item(2, 9).Value = "2023-01-11T05:32:11.915Z"
$ws.Cells.Item(2, 10).Value = "63be49dbc72b220018888298"
$ws.Cells.Item(3, 6).Value = 69.7
$ws.Cells.Item(3, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(3, 8).Value = "2023-01-11T05:32:11.926Z"
$ws.Cells.Item(3, 9).Value = "2023-01-11T05:32:11.926Z"
$ws.Cells.Item(3, 10).Value = "63be49dbc72b220018888299"
$ws.Cells.Item(4, 6).Value = 69.62
$ws.Cells.Item(4, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(4, 8).Value = "2023-01-11T05:32:11.933Z"
$ws.Cells.Item(4, 9).Value = "2023-01-11T05:32:11.933Z"
$ws.Cells.Item(4, 10).Value = "63be49dbc72b22001888829a"
$ws.Cells.Item(5, 6).Value = 70.12
$ws.Cells.Item(5, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(5, 8).Value = "2023-01-11T05:32:11.941Z"
$ws.Cells.Item(5, 9).Value = "2023-01-11T05:32:11.941Z"
$ws.Cells.Item(5, 10).Value = "63be49dbc72b22001888829b"
$ws.Cells.Item(6, 6).Value = 70.38
$ws.Cells.Item(6, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(6, 8).Value = "2023-01-11T05:32:11.948Z"
$ws.Cells.Item(6, 9).Value = "2023-01-11T05:32:11.948Z"
$ws.Cells.Item(6, 10).Value = "63be49dbc72b22001888829c"
$ws.Cells.Item(7, 6).Value = 72.75
$ws.Cells.Item(7, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(7, 8).Value = "2023-01-11T05:32:11.955Z"
$ws.Cells.Item(7, 9).Value = "2023-01-11T05:32:11.955Z"
$ws.Cells.Item(7, 10).Value = "63be49dbc72b22001888829d"
$ws.Cells.Item(8, 6).Value = 76.65
$ws.Cells.Item(8, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(8, 8).Value = "2023-01-11T05:32:11.962Z"
$ws.Cells.Item(8, 9).Value = "2023-01-11T05:32:11.962Z"
$ws.Cells.Item(8, 10).Value = "63be49dbc72b22001888829e"
$ws.Cells.Item(9, 6).Value = 76.65
$ws.Cells.Item(9, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(9, 8).Value = "2023-01-11T05:32:11.969Z"
$ws.Cells.Item(9, 9).Value = "2023-01-11T05:32:11.969Z"
$ws.Cells.Item(9, 10).Value = "63be49dbc72b22001888829f"
$ws.Cells.Item(10, 6).Value = 76.65
$ws.Cells.Item(10, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(10, 8).Value = "2023-01-11T05:32:11.976Z"
$ws.Cells.Item(10, 9).Value = "2023-01-11T05:32:11.976Z"
$ws.Cells.Item(10, 10).Value = "63be49dbc72b2200188882a0"
$ws.Cells.Item(11, 6).Value = 76.65
$ws.Cells.Item(11, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(11, 8).Value = "2023-01-11T05:32:11.982Z"
$ws.Cells.Item(11, 9).Value = "2023-01-11T05:32:11.982Z"
$ws.Cells.Item(11, 10).Value = "63be49dbc72b2200188882a1"
$ws.Cells.Item(12, 6).Value = 76.65
$ws.Cells.Item(12, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(12, 8).Value = "2023-01-11T05:32:11.989Z"
$ws.Cells.Item(12, 9).Value = "2023-01-11T05:32:11.989Z"
$ws.Cells.Item(12, 10).Value = "63be49dbc72b2200188882a2"
$ws.Cells.Item(13, 6).Value = 76.65
$ws.Cells.Item(13, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(13, 8).Value = "2023-01-11T05:32:11.997Z"
$ws.Cells.Item(13, 9).Value = "2023-01-11T05:32:11.997Z"
$ws.Cells.Item(13, 10).Value = "63be49dbc72b2200188882a3"
$ws.Cells.Item(14, 6).Value = 94.75
$ws.Cells.Item(14, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(14, 8).Value = "2023-01-11T05:32:12.004Z"
$ws.Cells.Item(14, 9).Value = "2023-01-11T05:32:12.004Z"
$ws.Cells.Item(14, 10).Value = "63be49dcc72b2200188882a4"
$ws.Cells.Item(15, 6).Value = 94.75
$ws.Cells.Item(15, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(15, 8).Value = "2023-01-11T05:32:12.015Z"
$ws.Cells.Item(15, 9).Value = "2023-01-11T05:32:12.015Z"
$ws.Cells.Item(15, 10).Value = "63be49dcc72b2200188882a5"
$ws.Cells.Item(16, 6).Value = 94.75
$ws.Cells.Item(16, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(16, 8).Value = "2023-01-11T05:32:12.022Z"
$ws.Cells.Item(16, 9).Value = "2023-01-11T05:32:12.022Z"
$ws.Cells.Item(16, 10).Value = "63be49dcc72b2200188882a6"
$ws.Cells.Item(17, 6).Value = 94.75
$ws.Cells.Item(17, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(17, 8).Value = "2023-01-11T05:32:12.029Z"
$ws.Cells.Item(17, 9).Value = "2023-01-11T05:32:12.029Z"
$ws.Cells.Item(17, 10).Value = "63be49dcc72b2200188882a7"
$ws.Cells.Item(18, 6).Value = 94.75
$ws.Cells.Item(18, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(18, 8).Value = "2023-01-11T05:32:12.036Z"
$ws.Cells.Item(18, 9).Value = "2023-01-11T05:32:12.036Z"
$ws.Cells.Item(18, 10).Value = "63be49dcc72b2200188882a8"
$ws.Cells.Item(19, 6).Value = 94.75
$ws.Cells.Item(19, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(19, 8).Value = "2023-01-11T05:32:12.043Z"
$ws.Cells.Item(19, 9).Value = "2023-01-11T05:32:12.043Z"
$ws.Cells.Item(19, 10).Value = "63be49dcc72b2200188882a9"
$ws.Cells.Item(20, 6).Value = 94.75
$ws.Cells.Item(20, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(20, 8).Value = "2023-01-11T05:32:12.050Z"
$ws.Cells.Item(20, 9).Value = "2023-01-11T05:32:12.050Z"
$ws.Cells.Item(20, 10).Value = "63be49dcc72b2200188882aa"
$ws.Cells.Item(21, 6).Value = 94.75
$ws.Cells.Item(21, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(21, 8).Value = "2023-01-11T05:32:12.058Z"
$ws.Cells.Item(21, 9).Value = "2023-01-11T05:32:12.058Z"
$ws.Cells.Item(21, 10).Value = "63be49dcc72b2200188882ab"
$ws.Cells.Item(22, 6).Value = 94.75
$ws.Cells.Item(22, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(22, 8).Value = "2023-01-11T05:32:12.065Z"
$ws.Cells.Item(22, 9).Value = "2023-01-11T05:32:12.065Z"
$ws.Cells.Item(22, 10).Value = "63be49dcc72b2200188882ac"
$ws.Cells.Item(23, 6).Value = 94.75
$ws.Cells.Item(23, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(23, 8).Value = "2023-01-11T05:32:12.072Z"
$ws.Cells.Item(23, 9).Value = "2023-01-11T05:32:12.072Z"
$ws.Cells.Item(23, 10).Value = "63be49dcc72b2200188882ad"
$ws.Cells.Item(24, 6).Value = 94.75
$ws.Cells.Item(24, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(24, 8).Value = "2023-01-11T05:32:12.079Z"
$ws.Cells.Item(24, 9).Value = "2023-01-11T05:32:12.079Z"
$ws.Cells.Item(24, 10).Value = "63be49dcc72b2200188882ae"
$ws.Cells.Item(25, 6).Value = 94.75
$ws.Cells.Item(25, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(25, 8).Value = "2023-01-11T05:32:12.086Z"
$ws.Cells.Item(25, 9).Value = "2023-01-11T05:32:12.086Z"
$ws.Cells.Item(25, 10).Value = "63be49dcc72b2200188882af"
$ws.Cells.Item(26, 6).Value = 107.5
$ws.Cells.Item(26, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(26, 8).Value = "2023-01-11T05:32:12.093Z"
$ws.Cells.Item(26, 9).Value = "2023-01-11T05:32:12.093Z"
$ws.Cells.Item(26, 10).Value = "63be49dcc72b2200188882b0"
$ws.Cells.Item(27, 6).Value = 107.5
$ws.Cells.Item(27, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(27, 8).Value = "2023-01-11T05:32:12.100Z"
$ws.Cells.Item(27, 9).Value = "2023-01-11T05:32:12.100Z"
$ws.Cells.Item(27, 10).Value = "63be49dcc72b2200188882b1"
$ws.Cells.Item(28, 6).Value = 107.5
$ws.Cells.Item(28, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(28, 8).Value = "2023-01-11T05:32:12.107Z"
$ws.Cells.Item(28, 9).Value = "2023-01-11T05:32:12.107Z"
$ws.Cells.Item(28, 10).Value = "63be49dcc72b2200188882b2"
$ws.Cells.Item(29, 6).Value = 107.5
$ws.Cells.Item(29, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(29, 8).Value = "2023-01-11T05:32:12.114Z"
$ws.Cells.Item(29, 9).Value = "2023-01-11T05:32:12.114Z"
$ws.Cells.Item(29, 10).Value = "63be49dcc72b2200188882b3"
$ws.Cells.Item(30, 6).Value = 107.5
$ws.Cells.Item(30, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(30, 8).Value = "2023-01-11T05:32:12.121Z"
$ws.Cells.Item(30, 9).Value = "2023-01-11T05:32:12.121Z"
$ws.Cells.Item(30, 10).Value = "63be49dcc72b2200188882b4"
$ws.Cells.Item(31, 6).Value = 107.5
$ws.Cells.Item(31, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(31, 8).Value = "2023-01-11T05:32:12.129Z"
$ws.Cells.Item(31, 9).Value = "2023-01-11T05:32:12.129Z"
$ws.Cells.Item(31, 10).Value = "63be49dcc72b2200188882b5"
$ws.Cells.Item(32, 6).Value = 107.5
$ws.Cells.Item(32, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(32, 8).Value = "2023-01-11T05:32:12.135Z"
$ws.Cells.Item(32, 9).Value = "2023-01-11T05:32:12.135Z"
$ws.Cells.Item(32, 10).Value = "63be49dcc72b2200188882b6"
$ws.Cells.Item(33, 6).Value = 107.5
$ws.Cells.Item(33, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(33, 8).Value = "2023-01-11T05:32:12.142Z"
$ws.Cells.Item(33, 9).Value = "2023-01-11T05:32:12.142Z"
$ws.Cells.Item(33, 10).Value = "63be49dcc72b2200188882b7"
$ws.Cells.Item(34, 6).Value = 107.5
$ws.Cells.Item(34, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(34, 8).Value = "2023-01-11T05:32:12.148Z"
$ws.Cells.Item(34, 9).Value = "2023-01-11T05:32:12.148Z"
$ws.Cells.Item(34, 10).Value = "63be49dcc72b2200188882b8"
$ws.Cells.Item(35, 6).Value = 107.5
$ws.Cells.Item(35, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(35, 8).Value = "2023-01-11T05:32:12.155Z"
$ws.Cells.Item(35, 9).Value = "2023-01-11T05:32:12.155Z"
$ws.Cells.Item(35, 10).Value = "63be49dcc72b2200188882b9"
$ws.Cells.Item(36, 6).Value = 107.5
$ws.Cells.Item(36, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(36, 8).Value = "2023-01-11T05:32:12.162Z"
$ws.Cells.Item(36, 9).Value = "2023-01-11T05:32:12.162Z"
$ws.Cells.Item(36, 10).Value = "63be49dcc72b2200188882ba"
$ws.Cells.Item(37, 6).Value = 107.5
$ws.Cells.Item(37, 7).Value = "2023-01-10T03:00:00.000Z"
$ws.Cells.Item(37, 8).Value = "2023-01-11T05:32:12.169Z"
$ws.Cells.Item(37, 9).Value = "2023-01-11T05:32:12.169Z"
$ws.Cells.Item(37, 10).Value = "63be49dcc72b2200188882bb"
